$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 88, which shifts the existing rows 88-180
# down to 89-181 (carrying their values/styles with them automatically).
$ws.Rows(88).Insert()

# Populate the newly inserted row 88 with the new data record.
$ws.Range("A88").Value = 8
$ws.Range("B88").Value = "Terminal La Palmera de La Serena"
$ws.Range("C88").Value = "Coquimbo"
$ws.Range("D88").Value = 44880
$ws.Range("E88").Value = 4
$ws.Range("F88").Value = 100112044
$ws.Range("G88").Value = "Perejil"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 2400
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = 1750
$ws.Range("N88").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O88").Value = "Provincia del Elquí"
$ws.Range("P88").Value = 1167
$ws.Range("Q88").Value = 1.5
$ws.Range("R88").Value = "Hortaliza"
